# Applies the "updates and add conclusion" revision to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row: insert a new "avx_opcount" column between M (shift_opcount)
# and the old N (global_index_count); SpeedUp stays in P.
# ---------------------------------------------------------------------
$ws.Range("O1").Value = $ws.Range("N1").Value2
$ws.Range("N1").Value = "avx_opcount"

# ---------------------------------------------------------------------
# Column widths: two new bestFit-ish columns (M / Q) picked up by the
# edit.
# ---------------------------------------------------------------------
$ws.Columns.Item(13).ColumnWidth = 8.79
$ws.Columns.Item(17).ColumnWidth = 8.79

# ---------------------------------------------------------------------
# Rows 26-30 ("Op Counts" block): drop the old E/F raw columns, add the
# new O (avx_opcount) column, and refresh the K..O op-count figures plus
# the D (cycle) baselines that feed the G/H/P formulas.
# ---------------------------------------------------------------------
$opRows = @(
  @{ Row=26; D=2359028;  K=9254482;  L=19958851;  M=3518120;  N=1344500; O=246394 },
  @{ Row=27; D=2604626;  K=13960000; L=25709593;  M=5851018;  N=1946000; O=293019 },
  @{ Row=28; D=3179679;  K=15843826; L=31034535;  M=6012434;  N=2334000; O=324016 },
  @{ Row=29; D=7252693;  K=35166241; L=61496918;  M=13041998; N=3701500; O=489502 },
  @{ Row=30; D=15988047; K=72765931; L=105889794; M=27779438; N=6971000; O=709055 }
)

foreach ($r in $opRows) {
  $row = $r.Row
  $ws.Range("E$row").ClearContents()
  $ws.Range("F$row").ClearContents()
  $ws.Range("D$row").Value2 = $r.D
  $ws.Range("K$row").Value2 = $r.K
  $ws.Range("L$row").Value2 = $r.L
  $ws.Range("M$row").Value2 = $r.M
  $ws.Range("N$row").Value2 = $r.N
  $ws.Range("O$row").Value2 = $r.O
  $ws.Range("G$row").Formula = "=K$row+L$row+M$row+N$row+O$row"
}

# G26 already carried style 5; the other rows (27-30) pick up the same
# highlighted style that G26/G30 used. Use a formats-only paste so the
# shared style index is reused instead of a near-duplicate xf being
# created.
$ws.Range("G26").Copy() | Out-Null
$ws.Range("G27:G30").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Rows 55-56 (benchmark block #1, label 101010 in A59 originally).
# ---------------------------------------------------------------------
$ws.Range("B55").Value2 = 1396200
$ws.Range("C55").Value2 = 1637724
$ws.Range("D55").Value2 = 1793954
$ws.Range("E55").Value2 = 4078117
$ws.Range("F55").Value2 = 8865370

$ws.Range("B56").Value2 = 1987830
$ws.Range("C56").Value2 = 2570753
$ws.Range("D56").Value2 = 3187398
$ws.Range("E56").Value2 = 6069194
$ws.Range("F56").Value2 = 13571951

# Row 59 label swaps from 101010 to 100000.
$ws.Range("A59").Value2 = 100000

$ws.Range("B60").Value2 = 975663
$ws.Range("C60").Value2 = 986119
$ws.Range("D60").Value2 = 1088488
$ws.Range("E60").Value2 = 2379518
$ws.Range("F60").Value2 = 5350095

$ws.Range("B61").Value2 = 1756744
$ws.Range("C61").Value2 = 2148763
$ws.Range("D61").Value2 = 2617051
$ws.Range("E61").Value2 = 5218790
$ws.Range("F61").Value2 = 11295713

# Row 65 label swaps from 100000 to 101010 (mirror of row 59).
$ws.Range("A65").Value2 = 101010

$ws.Range("B66").Value2 = 1038435
$ws.Range("C66").Value2 = 1148436
$ws.Range("D66").Value2 = 1315052
$ws.Range("E66").Value2 = 2800620
$ws.Range("F66").Value2 = 6149847

$ws.Range("B67").Value2 = 1812611
$ws.Range("C67").Value2 = 2232083
$ws.Range("D67").Value2 = 2838611
$ws.Range("E67").Value2 = 5673912
$ws.Range("F67").Value2 = 12057395

# ---------------------------------------------------------------------
# Selection / view: the author scrolled back to D30 (dropping the old
# topLeftCell/E57 selection).
# ---------------------------------------------------------------------
$ws.Range("D30").Select() | Out-Null
